$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
$ws.Range("D15").Value = "TIMESTAMP"
$ws.Range("D17").Value = "TIMESTAMP"
